$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 234.5
$ws.Range("I6").Value = 170
$ws.Range("K6").Value = 510
$ws.Range("M6").Value = -398

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 2771.0454
$ws.Range("I86").Value = 2537
$ws.Range("K86").Value = 2537
$ws.Range("M86").Value = -1414

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 2771.0454
$ws.Range("I89").Value = 2537
$ws.Range("K89").Value = 12685
$ws.Range("M89").Value = -7069

# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 457455.53
$ws.Range("J112").Value = 590995.5
$ws.Range("L112").Value = 1772986.5
$ws.Range("N112").Value = -1775202.5

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 6113.731
$ws.Range("I137").Value = 4861.773
$ws.Range("J137").Value = 12999.5
$ws.Range("K137").Value = 14585.319
$ws.Range("L137").Value = 38998.5
$ws.Range("M137").Value = -12035.319
$ws.Range("N137").Value = -44098.5

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 220467.12
$ws.Range("J138").Value = 326166.78
$ws.Range("L138").Value = 978500.3400000001
$ws.Range("N138").Value = -988780.3400000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 1986.5593
$ws.Range("I32").Value = 2083.18
$ws.Range("J32").Value = 1449.7778
$ws.Range("K32").Value = 2083.18
$ws.Range("L32").Value = 1449.7778
$ws.Range("M32").Value = -1796.18
$ws.Range("N32").Value = -2023.7778

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 8144
$ws.Range("I74").Value = 8101.3335
$ws.Range("J74").Value = 8176
$ws.Range("K74").Value = 8101.3335
$ws.Range("L74").Value = 8176
$ws.Range("M74").Value = -7227.3335
$ws.Range("N74").Value = -9924

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 8144
$ws.Range("I77").Value = 8101.3335
$ws.Range("J77").Value = 8176
$ws.Range("K77").Value = 40506.6675
$ws.Range("L77").Value = 40880
$ws.Range("M77").Value = -36138.6675
$ws.Range("N77").Value = -49616

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3531.1428
$ws.Range("I122").Value = 3311.125
$ws.Range("J122").Value = 3824.5
$ws.Range("K122").Value = 9933.375
$ws.Range("L122").Value = 11473.5
$ws.Range("M122").Value = -7483.375
$ws.Range("N122").Value = -16373.5

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 6079.85
$ws.Range("I132").Value = 3638.3635
$ws.Range("K132").Value = 10915.0905
$ws.Range("M132").Value = -8385.0905

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 6666.3335
$ws.Range("I99").Value = 4999.5
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 4999.5
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -3501.5
$ws.Range("N99").Value = -12996

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 4244.136
$ws.Range("I105").Value = 3342.4614
$ws.Range("K105").Value = 3342.4614
$ws.Range("M105").Value = -1595.4614

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1268.5294
$ws.Range("I16").Value = 1120.5385
$ws.Range("K16").Value = 1120.5385
$ws.Range("M16").Value = -833.5385000000001

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3497.5
$ws.Range("I31").Value = 2755.8333
$ws.Range("J31").Value = 3815.3572
$ws.Range("K31").Value = 2755.8333
$ws.Range("L31").Value = 3815.3572
$ws.Range("M31").Value = -2460.8333
$ws.Range("N31").Value = -4405.3572

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3497.5
$ws.Range("I34").Value = 2755.8333
$ws.Range("J34").Value = 3815.3572
$ws.Range("K34").Value = 2755.8333
$ws.Range("L34").Value = 3815.3572
$ws.Range("M34").Value = -2553.8333
$ws.Range("N34").Value = -4219.3572

# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 11426.818
$ws.Range("I86").Value = 7449.5
$ws.Range("J86").Value = 16199.6
$ws.Range("K86").Value = 7449.5
$ws.Range("L86").Value = 16199.6
$ws.Range("M86").Value = -6326.5
$ws.Range("N86").Value = -18445.6

# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 11426.818
$ws.Range("I89").Value = 7449.5
$ws.Range("J89").Value = 16199.6
$ws.Range("K89").Value = 37247.5
$ws.Range("L89").Value = 80998
$ws.Range("M89").Value = -31631.5
$ws.Range("N89").Value = -92230

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 405.24243
$ws.Range("J107").Value = 522.6923
$ws.Range("L107").Value = 522.6923
$ws.Range("N107").Value = -4362.6923

# Row 112 (Leve Item ID 25796)
$ws.Range("H112").Value = 90967.25
$ws.Range("J112").Value = 90967.25
$ws.Range("L112").Value = 90967.25
$ws.Range("N112").Value = -93921.25

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1268.5294
$ws.Range("I113").Value = 1120.5385
$ws.Range("K113").Value = 1120.5385
$ws.Range("M113").Value = 1049.4615

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 5557.1
$ws.Range("I134").Value = 5133.846
$ws.Range("J134").Value = 5880.7646
$ws.Range("K134").Value = 15401.538
$ws.Range("L134").Value = 17642.2938
$ws.Range("M134").Value = -12866.538
$ws.Range("N134").Value = -22712.2938

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 2653
$ws.Range("J68").Value = 2963.0908
$ws.Range("L68").Value = 8889.2724
$ws.Range("N68").Value = -10511.2724

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 2653
$ws.Range("J71").Value = 2963.0908
$ws.Range("L71").Value = 26667.8172
$ws.Range("N71").Value = -34779.8172

# Row 127 (Leve Item ID 38263)
$ws.Range("H127").Value = 1432.5
$ws.Range("J127").Value = 1432.5
$ws.Range("L127").Value = 4297.5
$ws.Range("N127").Value = -14217.5

$ws = $wb.Worksheets.Item("GSM")
# Row 21 (Leve Item ID 4430)
$ws.Range("H21").Value = 6007850
$ws.Range("J21").Value = 11675700
$ws.Range("L21").Value = 11675700
$ws.Range("N21").Value = -11676046

# Row 29 (Leve Item ID 4209)
$ws.Range("H29").Value = 28500
$ws.Range("J29").Value = 28500
$ws.Range("L29").Value = 28500
$ws.Range("N29").Value = -29080

# Row 30 (Leve Item ID 4430)
$ws.Range("H30").Value = 6007850
$ws.Range("J30").Value = 11675700
$ws.Range("L30").Value = 11675700
$ws.Range("N30").Value = -11675910

# Row 39 (Leve Item ID 18264)
$ws.Range("H39").Value = 49261
$ws.Range("J39").Value = 49261
$ws.Range("L39").Value = 49261
$ws.Range("N39").Value = -50325

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 22069.572
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 22069.572
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = 22069.572
$ws.Range("N80").Value = -24065.572

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 22069.572
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 22069.572
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = 110347.86
$ws.Range("N83").Value = -120331.86

# Row 116 (Leve Item ID 26120)
$ws.Range("H116").Value = 100000
$ws.Range("J116").Value = 100000
$ws.Range("L116").Value = 100000
$ws.Range("N116").Value = -109178

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 3695.8667
$ws.Range("I122").Value = 3254.1667
$ws.Range("J122").Value = 5462.6665
$ws.Range("K122").Value = 9762.500100000001
$ws.Range("L122").Value = 16387.9995
$ws.Range("M122").Value = -7312.500100000001
$ws.Range("N122").Value = -21287.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 9973.842000000001
$ws.Range("I16").Value = 15399.111
$ws.Range("K16").Value = 15399.111
$ws.Range("M16").Value = -15229.111

# Row 64 (Leve Item ID 10810)
$ws.Range("H64").Value = 39815.5
$ws.Range("J64").Value = 43778.6
$ws.Range("L64").Value = 43778.6
$ws.Range("N64").Value = -44228.6

# Row 67 (Leve Item ID 10810)
$ws.Range("H67").Value = 39815.5
$ws.Range("J67").Value = 43778.6
$ws.Range("L67").Value = 43778.6
$ws.Range("N67").Value = -45338.6

$ws = $wb.Worksheets.Item("WVR")
# Row 5 (Leve Item ID 3515)
$ws.Range("H5").Value = 1100000
$ws.Range("J5").Value = 2000000
$ws.Range("L5").Value = 2000000
$ws.Range("N5").Value = -2000224

# Row 64 (Leve Item ID 11036)
$ws.Range("H64").Value = 49975.25
$ws.Range("I64").Value = 20001
$ws.Range("K64").Value = 20001
$ws.Range("M64").Value = -19753

# Row 67 (Leve Item ID 11036)
$ws.Range("H67").Value = 49975.25
$ws.Range("I67").Value = 20001
$ws.Range("K67").Value = 20001
$ws.Range("M67").Value = -19143

# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 1262.375
$ws.Range("I100").Value = 931
$ws.Range("J100").Value = 1654
$ws.Range("K100").Value = 1862
$ws.Range("L100").Value = 3308
$ws.Range("M100").Value = -1321
$ws.Range("N100").Value = -4390

# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 1327.9286
$ws.Range("I113").Value = 1197.7
$ws.Range("K113").Value = 3593.1
$ws.Range("M113").Value = -1423.1

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 2107.394
$ws.Range("I122").Value = 1956.9
$ws.Range("J122").Value = 2338.923
$ws.Range("K122").Value = 5870.700000000001
$ws.Range("L122").Value = 7016.768999999999
$ws.Range("M122").Value = -3420.700000000001
$ws.Range("N122").Value = -11916.769
